# Generate Report for Handoff
# Updates the localization-status report after a new handoff: the zh-cn /
# de-de rows move from "Handed back: in sync with en-US" to "In Translation",
# their "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps advance,
# and the zh-cn / de-de "Error Detail" column picks up a stale-handback
# warning that used to be blank.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus      = "In Translation"
$newHoDate      = "2017-02-15 05:56:57"
$newZhHoDate    = "2017-02-15 05:56:39"
$errorDetail    = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a422ef894b3984df701e3a34ab28da698aa524d7/e2e/942e6272-aa00-41cc-9f69-e0acc0521251.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f0bcd3c205f0d44f6a373e02ac11ad6cc6d977b/e2e/942e6272-aa00-41cc-9f69-e0acc0521251.md."

# --- Overview sheet ---------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("G2").Value = $newHoDate

# --- zh-cn sheet --------------------------------------------------------
$ws2.Range("C2").Value = $newStatus
$ws2.Range("H2").Value = $newZhHoDate
$ws2.Range("R2").Value = $errorDetail

# --- de-de sheet --------------------------------------------------------
$ws3.Range("C2").Value = $newStatus
$ws3.Range("H2").Value = $newHoDate
$ws3.Range("R2").Value = $errorDetail

# --- Column width adjustments -------------------------------------------
# The Status columns shrink to fit the shorter "In Translation" text, and
# the Error Detail columns widen to fit the new long warning message.
$ws1.Columns.Item(5).ColumnWidth  = 12.5
$ws1.Columns.Item(6).ColumnWidth  = 12.5
$ws2.Columns.Item(3).ColumnWidth  = 12.5
$ws2.Columns.Item(18).ColumnWidth = 39.1666666666667
$ws3.Columns.Item(3).ColumnWidth  = 12.5
$ws3.Columns.Item(18).ColumnWidth = 39.1666666666667
